$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update processed fall values
$ws.Range("D8").Value = 396
$ws.Range("D9").Value = 377

# Update the active selection to D3
$ws.Range("D3").Select()
